$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 8 (pushes old rows 8-27 down to 10-29)
$ws.Rows("8:9").Insert()

# New row 8: "A Visa Backlog Abroad Is Taking a Toll Inside the U.S., Too"
$ws.Range("A8").Value = "A Visa Backlog Abroad Is Taking a Toll Inside the U.S., Too"
$ws.Range("B8").Value = "April 13, 2023"
$ws.Range("C8").Value = "The pileup has left visitors from places like Brazil, Colombia, India and Mexico waiting months, even a year or more, to visit family or do business in America."
$ws.Range("D8").Value = "https://static01.nyt.com/images/2023/04/12/travel/12visa-backlog-illo/12visa-backlog-illo-thumbWide.jpg?quality=75&auto=webp&disable=upscale"
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = $false

# New row 9: "New Bargain Lodges Offer Rustic Comfort in the Heart of Nature"
$ws.Range("A9").Value = "New Bargain Lodges Offer Rustic Comfort in the Heart of Nature"
$ws.Range("B9").Value = "April 12, 2023"
$ws.Range("C9").Value = "From the Catskill Mountains in New York to the deserts of Utah, new or expanding hotel companies are creating affordable accommodations for those who love the outdoors."
$ws.Range("D9").Value = "https://static01.nyt.com/images/2023/04/10/travel/oakImage-1681159827777/oakImage-1681159827777-thumbWide.jpg?quality=75&auto=webp&disable=upscale"
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = $false

# All phrase_count values (column E) become 0 for every data row
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 5).Value = 0
}
